$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change header text from "$/hour" to "Hourly rate"
$ws.Range("B1").Value = "Hourly rate"

# Update B2:B4 to use TEXT() formula formatting the hourly rate as a string like "$12.37"
$ws.Range("B2").Formula = '=TEXT(I2,"$00.00")'
$ws.Range("B3").Formula = '=TEXT(I3,"$00.00")'
$ws.Range("B4").Formula = '=TEXT(I4,"$00.00")'

# Update the active selection to B2
$ws.Range("B2").Select()
